$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (E1 / F1) ---
$ws.Range("E1").Value = "ENHANCE_ID"
$ws.Range("F1").Value = "Food sub-group"

# --- Data rows: ENHANCE_ID (col E) + Food sub-group (col F) ---
$data = @(
    @{Row=2; E=94850; F='Starchy staples'; Highlight=$false},
    @{Row=3; E=94879; F='Starchy staples'; Highlight=$false},
    @{Row=4; E=104945; F='Starchy staples'; Highlight=$false},
    @{Row=5; E=66362; F='Starchy staples'; Highlight=$false},
    @{Row=6; E=56375; F='Starchy staples'; Highlight=$false},
    @{Row=7; E=16440; F='Starchy staples'; Highlight=$false},
    @{Row=8; E=16446; F='Starchy staples'; Highlight=$false},
    @{Row=9; E=106464; F='Starchy staples'; Highlight=$false},
    @{Row=10; E=96496; F='Vegetables'; Highlight=$true},
    @{Row=11; E=65322; F='Legumes, nuts & seeds'; Highlight=$false},
    @{Row=12; E=55340; F='Legumes, nuts & seeds'; Highlight=$false},
    @{Row=13; E=55388; F='Legumes, nuts & seeds'; Highlight=$false},
    @{Row=14; E=55519; F='Legumes, nuts & seeds'; Highlight=$false},
    @{Row=15; E=15557; F='Legumes, nuts & seeds'; Highlight=$false},
    @{Row=16; E=55573; F='Legumes, nuts & seeds'; Highlight=$false},
    @{Row=17; E=56661; F='Vegetables'; Highlight=$false},
    @{Row=18; E=96689; F='Vegetables'; Highlight=$false},
    @{Row=19; E=96713; F='Vegetables'; Highlight=$false},
    @{Row=20; E=96715; F='Vegetables'; Highlight=$false},
    @{Row=21; E=56716; F='Vegetables'; Highlight=$false},
    @{Row=22; E=56746; F='Vegetables'; Highlight=$false},
    @{Row=23; E=56757; F='Vegetables'; Highlight=$false},
    @{Row=24; E=56787; F='Vegetables'; Highlight=$false},
    @{Row=25; E=56796; F='Vegetables'; Highlight=$false},
    @{Row=26; E=96797; F='Vegetables'; Highlight=$false},
    @{Row=27; E=16822; F='Vegetables'; Highlight=$false},
    @{Row=28; E=16832; F='Vegetables'; Highlight=$false},
    @{Row=29; E=16835; F='Vegetables'; Highlight=$false},
    @{Row=30; E=56839; F='Vegetables'; Highlight=$false},
    @{Row=31; E=1332; F='Vegetables'; Highlight=$false},
    @{Row=32; E=96849; F='Vegetables'; Highlight=$false},
    @{Row=33; E=16872; F='Vegetables'; Highlight=$false},
    @{Row=34; E=96876; F='Vegetables'; Highlight=$false},
    @{Row=35; E=16938; F='Vegetables'; Highlight=$false},
    @{Row=36; E=56957; F='Vegetables'; Highlight=$false},
    @{Row=37; E=96986; F='Vegetables'; Highlight=$false},
    @{Row=38; E=17000; F='Vegetables'; Highlight=$false},
    @{Row=39; E=97083; F='Vegetables'; Highlight=$false},
    @{Row=40; E=97084; F='Vegetables'; Highlight=$false},
    @{Row=41; E=17114; F='Vegetables'; Highlight=$false},
    @{Row=42; E=57123; F='Vegetables'; Highlight=$false},
    @{Row=43; E=57147; F='Vegetables'; Highlight=$false},
    @{Row=44; E=17166; F='Vegetables'; Highlight=$false},
    @{Row=45; E=57173; F='Vegetables'; Highlight=$false},
    @{Row=46; E=17258; F='Vegetables'; Highlight=$false},
    @{Row=47; E=44288; F='Fruits'; Highlight=$false},
    @{Row=48; E=54401; F='Fruits'; Highlight=$false},
    @{Row=49; E=14415; F='Fruits'; Highlight=$false},
    @{Row=50; E=44493; F='Fruits'; Highlight=$false},
    @{Row=51; E=54499; F='Fruits'; Highlight=$false},
    @{Row=52; E=14526; F='Fruits'; Highlight=$false},
    @{Row=53; E=14540; F='Fruits'; Highlight=$false},
    @{Row=54; E=14574; F='Fruits'; Highlight=$false},
    @{Row=55; E=74618; F='Fruits'; Highlight=$false},
    @{Row=56; E=17246; F='Fruits'; Highlight=$false},
    @{Row=57; E=55595; F='Animal source foods'; Highlight=$false},
    @{Row=58; E=15638; F='Animal source foods'; Highlight=$false},
    @{Row=59; E=15650; F='Animal source foods'; Highlight=$false},
    @{Row=60; E=15771; F='Animal source foods'; Highlight=$false},
    @{Row=61; E=55774; F='Animal source foods'; Highlight=$false},
    @{Row=62; E=95808; F='Animal source foods'; Highlight=$false},
    @{Row=63; E=15823; F='Animal source foods'; Highlight=$false},
    @{Row=64; E=95957; F='Animal source foods'; Highlight=$false},
    @{Row=65; E=95977; F='Animal source foods'; Highlight=$false},
    @{Row=66; E=96004; F='Animal source foods'; Highlight=$false},
    @{Row=67; E=55840; F='Animal source foods'; Highlight=$false},
    @{Row=68; E=56026; F='Animal source foods'; Highlight=$false},
    @{Row=69; E=13826; F='Animal source foods'; Highlight=$false},
    @{Row=70; E=103879; F='Animal source foods'; Highlight=$false},
    @{Row=71; E=13909; F='Animal source foods'; Highlight=$false},
    @{Row=72; E=13928; F='Animal source foods'; Highlight=$false},
    @{Row=73; E=74111; F='Animal source foods'; Highlight=$false},
    @{Row=74; E=54117; F='Animal source foods'; Highlight=$false},
    @{Row=75; E=14121; F='Animal source foods'; Highlight=$false},
    @{Row=76; E=54151; F='Animal source foods'; Highlight=$false},
    @{Row=77; E=94153; F='Animal source foods'; Highlight=$false},
    @{Row=78; E=54233; F='Animal source foods'; Highlight=$false},
    @{Row=79; E=54240; F='Animal source foods'; Highlight=$false},
    @{Row=80; E=76198; F='Animal source foods'; Highlight=$false},
    @{Row=81; E=26236; F='Animal source foods'; Highlight=$false},
    @{Row=82; E=13726; F='Animal source foods'; Highlight=$false},
    @{Row=83; E=13740; F='Animal source foods'; Highlight=$false},
    @{Row=84; E=56244; F='Oils & fats'; Highlight=$false},
    @{Row=85; E=96320; F='Oils & fats'; Highlight=$false},
    @{Row=86; E=56571; F=$null; Highlight=$true},
    @{Row=87; E=16574; F=$null; Highlight=$true},
    @{Row=88; E=55236; F=$null; Highlight=$true}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 5).Value = $item.E
    if ($item.F -ne $null) {
        $ws.Cells.Item($item.Row, 6).Value = $item.F
    }
    if ($item.Highlight) {
        $ws.Cells.Item($item.Row, 6).Interior.Color = 65535
    }
}

# --- Column E width adjustment ---
$ws.Columns("E:E").ColumnWidth = 13.86

# --- Clear the lingering cell selection left over from editing ---
$ws.Range("A1").Select()
